$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format while we assign values, to avoid Excel
# auto-converting numeric-looking strings (e.g. "593.10") into floating point
# numbers and losing formatting like trailing zeros / thousand separators.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '67.087.15'
$ws.Range('E2').Value = '  +5.97%  '
$ws.Range('D3').Value = '3.534.17'
$ws.Range('E3').Value = '  +4.22%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '593.10'
$ws.Range('E5').Value = '  +4.92%  '
$ws.Range('D6').Value = '170.41'
$ws.Range('E6').Value = '  +9.49%  '
$ws.Range('B7').Value = 'LidoStakedEther'
$ws.Range('C7').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D7').Value = '3.534.07'
$ws.Range('E7').Value = '  +4.20%  '
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('D9').Value = '0.583'
$ws.Range('E9').Value = '  +2.63%  '
$ws.Range('D10').Value = '7.27'
$ws.Range('E10').Value = '  +0.62%  '
$ws.Range('E11').Value = '  +6.61%  '
$ws.Range('D12').Value = '0.439'
$ws.Range('E12').Value = '  +4.47%  '
$ws.Range('D13').Value = '4.131.60'
$ws.Range('E13').Value = '  +3.64%  '
$ws.Range('E14').Value = '  +0.58%  '
$ws.Range('D15').Value = '28.36'
$ws.Range('E15').Value = '  +6.09%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '0.0000179'
$ws.Range('E16').Value = '  +5.73%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '67.021.55'
$ws.Range('E17').Value = '  +5.70%  '
$ws.Range('D18').Value = '3.538.79'
$ws.Range('E18').Value = '  +4.16%  '
$ws.Range('D19').Value = '6.31'
$ws.Range('E19').Value = '  +4.35%  '
$ws.Range('D20').Value = '14.11'
$ws.Range('E20').Value = '  +5.04%  '
$ws.Range('D21').Value = '392.94'
$ws.Range('E21').Value = '  +3.09%  '
$ws.Range('D22').Value = '7.98'
$ws.Range('E22').Value = '  +3.76%  '
$ws.Range('D23').Value = '73.40'
$ws.Range('E23').Value = '  +3.51%  '
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').Value = '0.0000126'
$ws.Range('E25').Value = '  +12.13%  '
$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D26').Value = '0.530'
$ws.Range('E26').Value = '  +3.66%  '
$ws.Range('D27').Value = '10.19'
$ws.Range('E27').Value = '  +5.88%  '
$ws.Range('D28').Value = '0.182'
$ws.Range('E28').Value = '  +2.45%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('E30').Value = '  +7.12%  '
$ws.Range('D31').Value = '1.48'
$ws.Range('E31').Value = '  +7.33%  '
$ws.Range('D32').Value = '2.07'
$ws.Range('E32').Value = '  +5.02%  '
$ws.Range('D33').Value = '23.61'
$ws.Range('E33').Value = '  +3.96%  '
$ws.Range('D34').Value = '7.42'
$ws.Range('E34').Value = '  +8.13%  '
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').Value = '1.59'
$ws.Range('E36').Value = '  +6.36%  '
$ws.Range('D37').Value = '161.61'
$ws.Range('E37').Value = '  +0.97%  '
$ws.Range('D38').Value = '0.912'
$ws.Range('E38').Value = '  +9.03%  '
$ws.Range('D39').Value = '1.95'
$ws.Range('E39').Value = '  +7.97%  '
$ws.Range('D40').Value = '0.0750'
$ws.Range('E40').Value = '  +5.62%  '
$ws.Range('D41').Value = '4.67'
$ws.Range('E41').Value = '  +8.40%  '
$ws.Range('D42').Value = '6.74'
$ws.Range('E42').Value = '  +6.33%  '
$ws.Range('D43').Value = '26.58'
$ws.Range('E43').Value = '  +3.09%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Value = '27.05'
$ws.Range('E44').Value = '  +6.79%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '2.822.44'
$ws.Range('E45').Value = '  +0.53%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value = '43.64'
$ws.Range('E46').Value = '  +1.76%  '
$ws.Range('D47').Value = '2.57'
$ws.Range('E47').Value = '  +12.87%  '
$ws.Range('D48').Value = '0.0316'
$ws.Range('E48').Value = '  +5.10%  '
$ws.Range('D49').Value = '356.30'
$ws.Range('E49').Value = '  +10.19%  '
$ws.Range('E50').Value = '  +7.31%  '
$ws.Range('D51').Value = '33.39'
$ws.Range('E51').Value = '  +12.36%  '

# Restore the original (default/no explicit number format) style so the
# cell styling matches the source workbook.
$priceRange.Style = "Normal"
